$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 2
$ws.Range("A2").Value = 2052980
$ws.Range("B2").Value = 16360

# Update the active selection to E4
$ws.Range("E4").Select()
